# LOM3230.xlsx — sync the "Objetivos / Programa resumido / Programa / Método /
# Critério / Norma de recuperação / Bibliografia" rows with the real course
# content (the sheet had been laid out with several rows shifted by one,
# reusing stray text from other fields) and insert the missing
# "Docentes responsáveis:" value row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row above the old row 13 ("Programa resumido:") so the
#    "Docentes responsáveis:" label (still sitting in A12) finally gets a
#    value row under it. This pushes the old rows 13-23 down to 14-24,
#    which already carries the correct row heights for every row below.
# ---------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# The inserted row copied A12's style into A13; the target layout has no
# cell in column A on row 13 at all, so drop it.
$ws.Range("A13").Clear()

# ---------------------------------------------------------------------
# 2) Give the new B13/C13 cells the same look as the rest of the table
#    (col B: normal weight, wrapped; col C: same but red) by copying the
#    formatting from the row right above (still col-B/col-C styled) and
#    then filling in the value.
# ---------------------------------------------------------------------
$ws.Range("B11:C11").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B13").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C13").Value = "5840726 - Cristina Bormio Nunes"

# ---------------------------------------------------------------------
# 3) Fix the mismatched value cells so every label lines up with its own
#    text instead of a neighbour's.
# ---------------------------------------------------------------------

# Objetivos: -> the actual Portuguese objective text (row 10)
$objetivosPt = "Apresentar as técnicas experimentais de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# Programa resumido: -> short summary (row 14, formerly "Semestral")
$resumidoPt = "Estudo das técnicas de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."
$ws.Range("B14").Value = $resumidoPt
$ws.Range("C14").Value = $resumidoPt

# Programa: -> full syllabus text (row 16, formerly the activation date)
$programaPt = "Propriedades elétricas: condutividade elétrica em metais puros, ligas metálicas e semicondutores,  e supercondutores; Efeito Hall; Lei de Ohm e dependência com a temperatura.Propriedades magnéticas: susceptibilidade magnética e magnetização c.c. Curvas de histerese de materiais magnéticos macios. Medidas de magnetostricção.Propriedades térmicas dos materiais:  expansão térmica."
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# Método: -> lab-method text (row 19, formerly the professor's name)
$metodoPt = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."
$ws.Range("B19").Value = $metodoPt
$ws.Range("C19").Value = $metodoPt

# Critério: -> grading criterion text (row 20, formerly the method text)
$criterioPt = "Média aritmética das notas dos relatórios de cada experimento"
$ws.Range("B20").Value = $criterioPt
$ws.Range("C20").Value = $criterioPt

# Norma de recuperação: -> makeup-exam text (row 21, formerly the criterion text)
$recuperacaoPt = "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("B21").Value = $recuperacaoPt
$ws.Range("C21").Value = $recuperacaoPt

# Bibliografia: -> reference list (row 22, formerly the makeup-exam text)
$bibliografiaPt = "HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000.RAYMOND A. SERWAY, CLEMENT J. MOSES, CURT A. MOYER. Modern Physics 3rd Edition,  Cengage Learning, Inc., 2005.SOLYMAR, L.; WALSH, D. Electrical Properties of Materials, Oxford University Press, 2009.NICOLA A. SPALDIN, Magnetic Materials, Fundamentals and Applications, SECOND EDITION, Cambridge University Press, 2011ROBERT, P. Electrical and Magnetic Properties of Materials, Artech House, 1998.SPEYER, R. Thermal Analysis of Materials, CRC Press, 1993."
$ws.Range("B22").Value = $bibliografiaPt
$ws.Range("C22").Value = $bibliografiaPt

# ---------------------------------------------------------------------
# 4) Column A was sized together with column B (min=1,max=2); split that
#    range so column A alone keeps its own width definition.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 59.83
